$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add Column C with the literal sum of columns A and B for rows 1-11
# (this replaces what used to be a separate "Add Column" activity/step)
for ($r = 1; $r -le 11; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value = $a + $b
}

# Update the selection to reflect the new annotated column C
$ws.Range("C1:C12").Select()
